$wb = $excel.ActiveWorkbook

# Updated "想去人数" (wish-to-go count) figures for both the "展览" sheet
# and the "全部类型" sheet, which mirror each other.
$updates = @{
    3  = 363
    4  = 1875
    7  = 192
    8  = 745
    10 = 358
    11 = 4459
    14 = 1243
    17 = 824
    19 = 445
    21 = 218
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
